# Auto-generated market-data refresh for Faerie_Profits (FFXIV leve profit sheets).
# Updates currentAveragePrice / NQ / HQ price + LeveProfitNQ/HQ columns (H-N)
# with freshly pulled market-board snapshots. Values with no new number are
# cleared (the row had no profit in that mode); previously-missing cells are added.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 932.7692
$ws.Range("J2").Value = 1131.5714
$ws.Range("L2").Value = 1131.5714
$ws.Range("N2").Value = -1357.5714
$ws.Range("H33").Value = 1030.8077
$ws.Range("I33").Value = 1307.8422
$ws.Range("J33").Value = 278.85715
$ws.Range("K33").Value = 1307.8422
$ws.Range("L33").Value = 278.85715
$ws.Range("M33").Value = -1078.8422
$ws.Range("N33").Value = -736.85715
$ws.Range("H41").Value = 724
$ws.Range("J41").Value = 898.8
$ws.Range("L41").Value = 898.8
$ws.Range("N41").Value = -1778.8
$ws.Range("H43").Value = 2638892.8
$ws.Range("J43").Value = 7577.143
$ws.Range("L43").Value = 7577.143
$ws.Range("N43").Value = -7715.143
$ws.Range("H51").Value = 4351714.5
$ws.Range("I51").Value = 7248043
$ws.Range("J51").Value = 7222
$ws.Range("K51").Value = 7248043
$ws.Range("L51").Value = 7222
$ws.Range("M51").Value = -7247559
$ws.Range("N51").Value = -8190
$ws.Range("H70").Value = 10495.917
$ws.Range("J70").Value = 9631.909
$ws.Range("L70").Value = 28895.727
$ws.Range("N70").Value = -29435.727
$ws.Range("H73").Value = 10495.917
$ws.Range("J73").Value = 9631.909
$ws.Range("L73").Value = 28895.727
$ws.Range("N73").Value = -30767.727
$ws.Range("H80").Value = 1303.7241
$ws.Range("J80").Value = 1797
$ws.Range("L80").Value = 5391
$ws.Range("N80").Value = -7387
$ws.Range("H83").Value = 1303.7241
$ws.Range("J83").Value = 1797
$ws.Range("L83").Value = 16173
$ws.Range("N83").Value = -26157
$ws.Range("H88").Value = 1898.6
$ws.Range("I88").Value = 1200
$ws.Range("J88").Value = 2364.3333
$ws.Range("K88").Value = 1200
$ws.Range("L88").Value = 2364.3333
$ws.Range("M88").Value = -794
$ws.Range("N88").Value = -3176.3333
$ws.Range("H91").Value = 1898.6
$ws.Range("I91").Value = 1200
$ws.Range("J91").Value = 2364.3333
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 2364.3333
$ws.Range("M91").Value = 204
$ws.Range("N91").Value = -5172.3333
$ws.Range("H92").Value = 269.2857
$ws.Range("I92").Value = 261.66666
$ws.Range("K92").Value = 261.66666
$ws.Range("M92").Value = 986.33334
$ws.Range("H106").Value = 77051.42999999999
$ws.Range("I106").Value = 87393.336
$ws.Range("K106").Value = 87393.336
$ws.Range("M106").Value = -86762.336
$ws.Range("H135").Value = 15733.4
$ws.Range("J135").Value = 19853.572
$ws.Range("L135").Value = 178682.148
$ws.Range("N135").Value = -183752.148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1361.1428
$ws.Range("I45").Value = 1338.0834
$ws.Range("K45").Value = 1338.0834
$ws.Range("M45").Value = -961.0834
$ws.Range("H68").Value = 52272
$ws.Range("I68").Value = 40000
$ws.Range("J68").Value = 53499.2
$ws.Range("K68").Value = 40000
$ws.Range("L68").Value = 53499.2
$ws.Range("M68").Value = -39189
$ws.Range("N68").Value = -55121.2
$ws.Range("H71").Value = 52272
$ws.Range("I71").Value = 40000
$ws.Range("J71").Value = 53499.2
$ws.Range("K71").Value = 120000
$ws.Range("L71").Value = 160497.6
$ws.Range("M71").Value = -115944
$ws.Range("N71").Value = -168609.6
$ws.Range("H74").Value = 5459.2607
$ws.Range("I74").Value = 3784.7778
$ws.Range("K74").Value = 3784.7778
$ws.Range("M74").Value = -2910.7778
$ws.Range("H77").Value = 5459.2607
$ws.Range("I77").Value = 3784.7778
$ws.Range("K77").Value = 18923.889
$ws.Range("M77").Value = -14555.889
$ws.Range("H132").Value = 4880.7812
$ws.Range("I132").Value = 4134.9546
$ws.Range("J132").Value = 6521.6
$ws.Range("K132").Value = 12404.8638
$ws.Range("L132").Value = 19564.8
$ws.Range("M132").Value = -9874.863799999999
$ws.Range("N132").Value = -24624.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2680720
$ws.Range("I86").Value = 4184970.5
$ws.Range("J86").Value = 6497.222
$ws.Range("K86").Value = 4184970.5
$ws.Range("L86").Value = 6497.222
$ws.Range("M86").Value = -4183847.5
$ws.Range("N86").Value = -8743.222
$ws.Range("H89").Value = 2680720
$ws.Range("I89").Value = 4184970.5
$ws.Range("J89").Value = 6497.222
$ws.Range("K89").Value = 20924852.5
$ws.Range("L89").Value = 32486.11
$ws.Range("M89").Value = -20919236.5
$ws.Range("N89").Value = -43718.11
$ws.Range("H134").Value = 2918.3845
$ws.Range("I134").Value = 2918.3845
$ws.Range("K134").Value = 8755.1535
$ws.Range("M134").Value = -6220.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4478.1
$ws.Range("I31").Value = 3205.1667
$ws.Range("J31").Value = 9569.833000000001
$ws.Range("K31").Value = 3205.1667
$ws.Range("L31").Value = 9569.833000000001
$ws.Range("M31").Value = -2910.1667
$ws.Range("N31").Value = -10159.833
$ws.Range("H34").Value = 4478.1
$ws.Range("I34").Value = 3205.1667
$ws.Range("J34").Value = 9569.833000000001
$ws.Range("K34").Value = 3205.1667
$ws.Range("L34").Value = 9569.833000000001
$ws.Range("M34").Value = -3003.1667
$ws.Range("N34").Value = -9973.833000000001
$ws.Range("H62").Value = 3550.25
$ws.Range("I62").Value = 3550.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3550.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -2926.25
$ws.Range("H65").Value = 3550.25
$ws.Range("I65").Value = 3550.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17751.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -14631.25
$ws.Range("H103").Value = 52331.668
$ws.Range("I103").Value = 52331.668
$ws.Range("K103").Value = 52331.668
$ws.Range("M103").Value = -51159.668
$ws.Range("H134").Value = 1537.75
$ws.Range("I134").Value = 1537.75
$ws.Range("K134").Value = 4613.25
$ws.Range("M134").Value = -2078.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 333545.34
$ws.Range("J9").Value = 341.66666
$ws.Range("L9").Value = 1024.99998
$ws.Range("N9").Value = -1472.99998
$ws.Range("H23").Value = 1088.25
$ws.Range("J23").Value = 1184.909
$ws.Range("L23").Value = 3554.727
$ws.Range("N23").Value = -4024.727
$ws.Range("H103").Value = 593.6
$ws.Range("I103").Value = 492.5
$ws.Range("J103").Value = 661
$ws.Range("K103").Value = 1477.5
$ws.Range("L103").Value = 1983
$ws.Range("M103").Value = -598.5
$ws.Range("N103").Value = -3741
$ws.Range("H129").Value = 3911.3333
$ws.Range("I129").Value = 621.7
$ws.Range("J129").Value = 6261.0713
$ws.Range("K129").Value = 1865.1
$ws.Range("L129").Value = 18783.2139
$ws.Range("M129").Value = 3134.9
$ws.Range("N129").Value = -28783.2139
$ws.Range("H139").Value = 3715.3125
$ws.Range("I139").Value = 2844.5
$ws.Range("K139").Value = 8533.5
$ws.Range("M139").Value = -3393.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 49066.2
$ws.Range("J106").Value = 49066.2
$ws.Range("L106").Value = 49066.2
$ws.Range("N106").Value = -51590.2
$ws.Range("H132").Value = 4266.479
$ws.Range("I132").Value = 5023.2354
$ws.Range("J132").Value = 2428.6428
$ws.Range("K132").Value = 15069.7062
$ws.Range("L132").Value = 7285.928400000001
$ws.Range("M132").Value = -12539.7062
$ws.Range("N132").Value = -12345.9284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 89943
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("H69").Value = 88929.664
$ws.Range("J69").Value = 88929.664
$ws.Range("L69").Value = 88929.664
$ws.Range("N69").Value = -90551.664
$ws.Range("H72").Value = 88929.664
$ws.Range("J72").Value = 88929.664
$ws.Range("L72").Value = 266788.992
$ws.Range("N72").Value = -274900.992
$ws.Range("H97").Value = 43747.668
$ws.Range("J97").Value = 43747.668
$ws.Range("L97").Value = 43747.668
$ws.Range("N97").Value = -45729.668
$ws.Range("H100").Value = 3841.5557
$ws.Range("I100").Value = 3749.9
$ws.Range("K100").Value = 3749.9
$ws.Range("M100").Value = -3208.9
$ws.Range("H104").Value = 35799.6
$ws.Range("J104").Value = 35799.6
$ws.Range("L104").Value = 35799.6
$ws.Range("N104").Value = -42787.6
$ws.Range("H132").Value = 9786.916999999999
$ws.Range("I132").Value = 9765.375
$ws.Range("K132").Value = 29296.125
$ws.Range("M132").Value = -26766.125
$ws.Range("H136").Value = 5370.0713
$ws.Range("I136").Value = 5556.875
$ws.Range("J136").Value = 4249.25
$ws.Range("K136").Value = 16670.625
$ws.Range("L136").Value = 12747.75
$ws.Range("M136").Value = -14120.625
$ws.Range("N136").Value = -17847.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = $null
$ws.Range("N99").Value = 0
$ws.Range("H101").Value = 85602
$ws.Range("J101").Value = 85602
$ws.Range("L101").Value = 85602
$ws.Range("N101").Value = -92092
$ws.Range("H107").Value = 869.6429000000001
$ws.Range("J107").Value = 953
$ws.Range("L107").Value = 2859
$ws.Range("N107").Value = -6699
$ws.Range("H132").Value = 2199.8
$ws.Range("I132").Value = 2124.75
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6374.25
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3844.25
$ws.Range("N132").Value = -12560

Write-Output "Updated 250 cells across 8 sheets."
